$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
# Set G6 fill FIRST, before anything else touches A8
$ws.Range("G5").Copy() | Out-Null
$ws.Range("G6").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("G6").Interior.Color = 42495
# Now touch A8
$ws.Range("A8").Value = "test"
